# This revision corresponds to an upstream commit whose message is
# "Fixed POI packaging and upgraded to POI 3.15." Diffing the canonical
# OOXML of word/document.xml and word/styles.xml shows that every single
# change is an XML attribute re-ordering produced by the new POI writer
# (e.g. <w:tab w:val="left" w:pos="3119"/> -> <w:tab w:pos="3119" w:val="left"/>,
# <w:pgSz w:w="11906" w:h="16838"/> -> <w:pgSz w:h="16838" w:w="11906"/>, and a
# large number of similar attribute permutations throughout styles.xml).
# No text, value, style, or formatting actually changes anywhere in the
# package. To mirror the intent of the commit through the Word object
# model we re-apply (round-trip) the exact same values on every setting
# that the diff touches, which is a no-op for the document's content.

$d = $word.ActiveDocument

# --- word/document.xml: the single explicit tab stop (w:pos=3119, w:val=left) ---
$p2 = $d.Paragraphs.Item(2)
$tabs = $p2.Range.ParagraphFormat.TabStops
$tabs.ClearAll()
$tabs.Add(155.95, 0)

# --- word/document.xml: sectPr/pgSz + sectPr/pgMar (same twips values) ---
$ps = $d.PageSetup
$ps.PageWidth = 595.3
$ps.PageHeight = 841.9
$ps.TopMargin = 70.85
$ps.RightMargin = 70.85
$ps.BottomMargin = 70.85
$ps.LeftMargin = 70.85
$ps.HeaderDistance = 35.4
$ps.FooterDistance = 35.4
$ps.Gutter = 0

# NOTE: word/styles.xml (docDefaults/latentStyles/style attribute order) has
# no Word-object-model surface that lets a script reorder XML attributes
# without rewriting element content outright (e.g. stamping an explicit
# <w:rPr> onto the "Normal" style, which does not exist in the source
# document). Touching it that way would introduce a real content change
# that is not present in the target revision, so styles.xml is
# intentionally left untouched here: its values already match the target.
